$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "Semestre ideal:" value EQD-8,EQN-11 -> EQD-8,EQN-10
$ws.Range("B9").Value = "EQD-8,EQN-10"
$ws.Range("C9").Value = "EQD-8,EQN-10"

# 2. Requisitos: replace "LOQ4086 - Operacoes Unitarias II" entry with "LOQ4002 - Reatores Quimicos"
$ws.Range("B25").Value = "LOQ4002 -  Reatores Quimicos  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOQ4002 -  Reatores Quimicos  (Requisito fraco)`n"

# 3. Add a new requirement row 26: "LOQ4057 - Operacoes Unitarias III"
$ws.Range("B25:C25").Copy()
$ws.Range("B26:C26").PasteSpecial(-4122)
$ws.Range("B26").Value = "LOQ4057 -  Operações Unitárias III  (Requisito fraco)`n"
$ws.Range("C26").Value = "LOQ4057 -  Operações Unitárias III  (Requisito fraco)`n"
$ws.Rows.Item(26).RowHeight = 30
